$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 34496704
$ws.Range("I62").Value = 55563692
$ws.Range("K62").Value = 55563692
$ws.Range("M62").Value = -55563068

$ws.Range("H65").Value = 34496704
$ws.Range("I65").Value = 55563692
$ws.Range("K65").Value = 277818460
$ws.Range("M65").Value = -277815340

$ws.Range("H70").Value = 3138.8572
$ws.Range("I70").Value = 3303
$ws.Range("J70").Value = 3037.8462
$ws.Range("K70").Value = 9909
$ws.Range("L70").Value = 9113.5386
$ws.Range("M70").Value = -9639
$ws.Range("N70").Value = -9653.5386

$ws.Range("H73").Value = 3138.8572
$ws.Range("I73").Value = 3303
$ws.Range("J73").Value = 3037.8462
$ws.Range("K73").Value = 9909
$ws.Range("L73").Value = 9113.5386
$ws.Range("M73").Value = -8973
$ws.Range("N73").Value = -10985.5386

$ws.Range("H135").Value = 3249.476
$ws.Range("I135").Value = 3546.611
$ws.Range("K135").Value = 31919.499
$ws.Range("M135").Value = -29384.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2226.9106
$ws.Range("I32").Value = 2216.509
$ws.Range("K32").Value = 2216.509
$ws.Range("M32").Value = -1929.509

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H61").Value = 4797.032
$ws.Range("I61").Value = 5051.8
$ws.Range("K61").Value = 5051.8
$ws.Range("M61").Value = -4839.8

$ws.Range("H102").Value = 7706.591
$ws.Range("I102").Value = 4292.143
$ws.Range("K102").Value = 4292.143
$ws.Range("M102").Value = -2670.143

$ws.Range("H136").Value = 4797.032
$ws.Range("I136").Value = 5051.8
$ws.Range("K136").Value = 15155.4
$ws.Range("M136").Value = -12605.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1966.4615
$ws.Range("I20").Value = 1216.2778
$ws.Range("K20").Value = 1216.2778
$ws.Range("M20").Value = -969.2778000000001

$ws.Range("H82").Value = 55640.25
$ws.Range("J82").Value = 96282.75
$ws.Range("L82").Value = 96282.75
$ws.Range("N82").Value = -97048.75

$ws.Range("H85").Value = 55640.25
$ws.Range("J85").Value = 96282.75
$ws.Range("L85").Value = 96282.75
$ws.Range("N85").Value = -98934.75

$ws.Range("H134").Value = 5223.8047
$ws.Range("I134").Value = 5397.914
$ws.Range("K134").Value = 16193.742
$ws.Range("M134").Value = -13658.742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 745.56525
$ws.Range("I7").Value = 1090.4667
$ws.Range("J7").Value = 98.875
$ws.Range("K7").Value = 1090.4667
$ws.Range("L7").Value = 98.875
$ws.Range("M7").Value = -977.4666999999999
$ws.Range("N7").Value = -324.875

$ws.Range("H22").Value = 465.94736
$ws.Range("J22").Value = 473.63635
$ws.Range("L22").Value = 473.63635
$ws.Range("N22").Value = -1173.63635

$ws.Range("H124").Value = 49996
$ws.Range("J124").Value = 49996
$ws.Range("L124").Value = 49996
$ws.Range("N124").Value = -54906

$ws.Range("H132").Value = 33315.35
$ws.Range("I132").Value = 10397.5625
$ws.Range("K132").Value = 31192.6875
$ws.Range("M132").Value = -28662.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3715203.8
$ws.Range("I4").Value = 3715203.8
$ws.Range("K4").Value = 11145611.4
$ws.Range("M4").Value = -11145499.4

$ws.Range("H35").Value = 176
$ws.Range("I35").Value = 176
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 528
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -240
$ws.Range("N35").ClearContents()

$ws.Range("H41").Value = 4809.7617
$ws.Range("I41").Value = 8782.637000000001
$ws.Range("J41").Value = 439.6
$ws.Range("K41").Value = 26347.911
$ws.Range("L41").Value = 1318.8
$ws.Range("M41").Value = -26009.911
$ws.Range("N41").Value = -1994.8

$ws.Range("H52").Value = 41000
$ws.Range("J52").Value = 41000
$ws.Range("L52").Value = 123000
$ws.Range("N52").Value = -123532

$ws.Range("H55").Value = 7146.7144
$ws.Range("I55").Value = 917.8
$ws.Range("J55").Value = 10607.223
$ws.Range("K55").Value = 2753.4
$ws.Range("L55").Value = 31821.669
$ws.Range("M55").Value = -2576.4
$ws.Range("N55").Value = -32175.669

$ws.Range("H107").Value = 2774.875
$ws.Range("J107").Value = 2899.8333
$ws.Range("L107").Value = 8699.499899999999
$ws.Range("N107").Value = -12539.4999

$ws.Range("H117").Value = 1448
$ws.Range("I117").Value = 1448
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 4344
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -902
$ws.Range("N117").ClearContents()

$ws.Range("H128").Value = 349895.25
$ws.Range("I128").Value = 349895.25
$ws.Range("K128").Value = 1049685.75
$ws.Range("M128").Value = -1044705.75

$ws.Range("H129").Value = 27779224
$ws.Range("J129").Value = 55557892
$ws.Range("L129").Value = 166673676
$ws.Range("N129").Value = -166683676

$ws.Range("H131").Value = 23257288
$ws.Range("I131").Value = 142858370
$ws.Range("J131").Value = 1521.8334
$ws.Range("K131").Value = 428575110
$ws.Range("L131").Value = 4565.5002
$ws.Range("M131").Value = -428570070
$ws.Range("N131").Value = -14645.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 45331
$ws.Range("J93").Value = 45331
$ws.Range("L93").Value = 45331
$ws.Range("N93").Value = -49075

$ws.Range("H122").Value = 8635.486999999999
$ws.Range("I122").Value = 9894.532999999999
$ws.Range("J122").Value = 4438.6665
$ws.Range("K122").Value = 29683.599
$ws.Range("L122").Value = 13315.9995
$ws.Range("M122").Value = -27233.599
$ws.Range("N122").Value = -18215.9995

$ws.Range("H126").Value = 21620.736
$ws.Range("J126").Value = 18176.273
$ws.Range("L126").Value = 54528.819
$ws.Range("N126").Value = -59468.819

$ws.Range("H132").Value = 4434.3335
$ws.Range("I132").Value = 3260.9
$ws.Range("J132").Value = 6781.2
$ws.Range("K132").Value = 9782.700000000001
$ws.Range("L132").Value = 20343.6
$ws.Range("M132").Value = -7252.700000000001
$ws.Range("N132").Value = -25403.6

$ws.Range("H133").Value = 89998
$ws.Range("J133").Value = 89998
$ws.Range("L133").Value = 89998
$ws.Range("N133").Value = -100118

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22843.5
$ws.Range("I7").Value = 24585.96
$ws.Range("K7").Value = 24585.96
$ws.Range("M7").Value = -24473.96

$ws.Range("H68").Value = 4992.7856
$ws.Range("J68").Value = 6085.143
$ws.Range("L68").Value = 6085.143
$ws.Range("N68").Value = -7583.143

$ws.Range("H71").Value = 4992.7856
$ws.Range("J71").Value = 6085.143
$ws.Range("L71").Value = 30425.715
$ws.Range("N71").Value = -37913.715

$ws.Range("H93").Value = 5036.25
$ws.Range("I93").Value = 5978.6
$ws.Range("K93").Value = 5978.6
$ws.Range("M93").Value = -4730.6

$ws.Range("H126").Value = 22843.5
$ws.Range("I126").Value = 24585.96
$ws.Range("K126").Value = 73757.88
$ws.Range("M126").Value = -71287.88

$ws.Range("H132").Value = 1003462.44
$ws.Range("I132").Value = 2999948.5
$ws.Range("K132").Value = 8999845.5
$ws.Range("M132").Value = -8997315.5

$ws.Range("H136").Value = 7967.923
$ws.Range("I136").Value = 3099.1667
$ws.Range("J136").Value = 12141.143
$ws.Range("K136").Value = 9297.500100000001
$ws.Range("L136").Value = 36423.429
$ws.Range("M136").Value = -6747.500100000001
$ws.Range("N136").Value = -41523.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 23398.8
$ws.Range("I51").Value = 14249.75
$ws.Range("K51").Value = 14249.75
$ws.Range("M51").Value = -13739.75

$ws.Range("H122").Value = 2322.4285
$ws.Range("I122").Value = 1754.862
$ws.Range("J122").Value = 5065.6665
$ws.Range("K122").Value = 5264.586
$ws.Range("L122").Value = 15196.9995
$ws.Range("M122").Value = -2814.586
$ws.Range("N122").Value = -20096.9995

$ws.Range("H126").Value = 25088.475
$ws.Range("I126").Value = 31049.715
$ws.Range("J126").Value = 8397
$ws.Range("K126").Value = 93149.145
$ws.Range("L126").Value = 25191
$ws.Range("M126").Value = -90679.145
$ws.Range("N126").Value = -30131

$ws.Range("H132").Value = 24173.154
$ws.Range("I132").Value = 27024.428
$ws.Range("J132").Value = 12197.8
$ws.Range("K132").Value = 81073.284
$ws.Range("L132").Value = 36593.39999999999
$ws.Range("M132").Value = -78543.284
$ws.Range("N132").Value = -41653.39999999999

$ws.Range("H136").Value = 4459.7334
$ws.Range("I136").Value = 3190.7
$ws.Range("K136").Value = 9572.099999999999
$ws.Range("M136").Value = -7022.099999999999
